# Apply updated crypto price/volume data, and re-sort three swapped-rank rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.275.33"
$ws.Range("E2").Value = "  +4.24%  "
$ws.Range("D3").Value = "1.810.98"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.67%  "
$ws.Range("D5").Value = "'339.54"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").Value = "'0.9990"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").Value = "'0.3927"
$ws.Range("E7").Value = "  +4.44%  "
$ws.Range("D8").Value = "'0.3500"
$ws.Range("E8").Value = "  +3.16%  "
$ws.Range("D9").Value = "'48.60"
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("D10").Value = "'1.183"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").Value = "'0.07558"
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("D12").Value = "'0.9996"
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").Value = "'22.12"
$ws.Range("E13").Value = "  +3.68%  "
$ws.Range("D14").Value = "'6.538"
$ws.Range("E14").Value = "  +2.86%  "
$ws.Range("D15").Value = "1.813.28"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("D16").Value = "'7.180"
$ws.Range("E16").Value = "  +2.92%  "
$ws.Range("D17").Value = "'0.00001106"
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("D18").Value = "'0.06722"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").Value = "'85.36"
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").Value = "'0.9992"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").Value = "'17.81"
$ws.Range("E21").Value = "  +4.34%  "
$ws.Range("D22").Value = "'6.590"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").Value = "28.284.23"
$ws.Range("E23").Value = "  +4.25%  "
$ws.Range("D24").Value = "'12.48"
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("D25").Value = "'2.403"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("D26").Value = "'1.483"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.531"
$ws.Range("E27").Value = "  +1.69%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'21.42"
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("D29").Value = "'155.00"
$ws.Range("E29").Value = "  +2.43%  "
$ws.Range("D30").Value = "2.020.02"
$ws.Range("E30").Value = "  +1.94%  "
$ws.Range("D31").Value = "'136.20"
$ws.Range("E31").Value = "  +3.07%  "
$ws.Range("D32").Value = "'6.388"
$ws.Range("E32").Value = "  +7.89%  "
$ws.Range("D33").Value = "'4.021"
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("D34").Value = "'0.08863"
$ws.Range("E34").Value = "  +2.89%  "
$ws.Range("D35").Value = "'13.21"
$ws.Range("E35").Value = "  +2.65%  "
$ws.Range("D36").Value = "'0.02476"
$ws.Range("E36").Value = "  +7.06%  "
$ws.Range("D37").Value = "'5.499"
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").Value = "'0.6957"
$ws.Range("E38").Value = "  +3.11%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06575"
$ws.Range("E39").Value = "  +4.78%  "
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("D41").Value = "'0.2234"
$ws.Range("E41").Value = "  +3.23%  "
$ws.Range("D42").Value = "'1.268"
$ws.Range("E42").Value = "  +2.72%  "
$ws.Range("E43").Value = "  -1.23%  "
$ws.Range("D44").Value = "'14.72"
$ws.Range("E44").Value = "  +2.44%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6445"
$ws.Range("E45").Value = "  +2.70%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "'0.9986"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("D47").Value = "'3.875"
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("D48").Value = "'2.169"
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("D49").Value = "'131.99"
$ws.Range("E49").Value = "  +3.03%  "
$ws.Range("D50").Value = "'0.07248"
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("D51").Value = "'80.48"
$ws.Range("E51").Value = "  +2.73%  "
